$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.576.78'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.347.00'
$ws.Range('E3').Value = '  -2.94%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.33'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.00'
$ws.Range('E6').Value = '  -3.83%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.59'
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('E12').Value = '  -4.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.92'
$ws.Range('E13').Value = '  -5.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.762.96'
$ws.Range('E14').Value = '  -2.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.564.41'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.345.06'
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '318.77'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.60'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '63.95'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('E24').Value = '  -3.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.32'
$ws.Range('E26').Value = '  -3.35%  '
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.34'
$ws.Range('E29').Value = '  +1.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0746'
$ws.Range('E30').Value = '  -3.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.94'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('E32').Value = '  +6.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.399'
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.54'
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '313.21'
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '143.98'
$ws.Range('E42').Value = '  +2.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.47'
$ws.Range('E43').Value = '  -5.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0957'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('E45').Value = '  -2.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.84'
$ws.Range('E46').Value = '  -3.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.562'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  -0.35%  '
